# Updates the "Price" (column D) and "Volume(1h)" (column E) figures for the
# cryptocurrency rows on the active sheet, matching the latest scraped data.
#
# Notes:
#  - Column E values are percentage-looking strings padded with spaces
#    (e.g. "  -0.87%  "), so they are written as-is via Value2.
#  - Column D values are numeric-looking text (e.g. "305.33") that must stay
#    as plain text (not be converted into floating point numbers by Excel).
#    We force text interpretation with a leading apostrophe and then reset
#    the cell style back to "Normal" so no extra number formatting/style is
#    left applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'42.654.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  -0.87%  "
$ws.Range("D3").Value2 = "'2.280.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  -0.87%  "
$ws.Range("E4").Value2 = "  +0.01%  "
$ws.Range("D5").Value2 = "'305.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +1.75%  "
$ws.Range("D6").Value2 = "'96.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  -1.44%  "
$ws.Range("D7").Value2 = "'0.507"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  -2.78%  "
$ws.Range("E8").Value2 = "  +0.01%  "
$ws.Range("E9").Value2 = "  -3.37%  "
$ws.Range("D10").Value2 = "'35.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  -2.07%  "
$ws.Range("D11").Value2 = "'0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  -0.23%  "
$ws.Range("D12").Value2 = "'18.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  +2.64%  "
$ws.Range("E13").Value2 = "  +1.05%  "
$ws.Range("D14").Value2 = "'6.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -2.22%  "
$ws.Range("D15").Value2 = "'2.634.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  -0.87%  "
$ws.Range("D16").Value2 = "'2.271.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  -1.25%  "
$ws.Range("D17").Value2 = "'0.779"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -1.30%  "
$ws.Range("D18").Value2 = "'42.592.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -0.74%  "
$ws.Range("D19").Value2 = "'12.94"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -0.26%  "
$ws.Range("E20").Value2 = "  -1.73%  "
$ws.Range("D22").Value2 = "'67.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  -1.67%  "
$ws.Range("D23").Value2 = "'235.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  -0.79%  "
$ws.Range("E24").Value2 = "  -3.03%  "
$ws.Range("E25").Value2 = "  +1.54%  "
$ws.Range("E26").Value2 = "  +0.14%  "
$ws.Range("D27").Value2 = "'4.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +0.22%  "
$ws.Range("D28").Value2 = "'25.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +0.62%  "
$ws.Range("D29").Value2 = "'165.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +1.40%  "
$ws.Range("E30").Value2 = "  +0.06%  "
$ws.Range("D32").Value2 = "'33.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  -0.21%  "
$ws.Range("E33").Value2 = "  +0.09%  "
$ws.Range("E34").Value2 = "  -0.30%  "
$ws.Range("E35").Value2 = "  -3.10%  "
$ws.Range("D36").Value2 = "'17.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  -3.53%  "
$ws.Range("E37").Value2 = "  -0.98%  "
$ws.Range("D38").Value2 = "'0.0689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -1.13%  "
$ws.Range("E39").Value2 = "  -1.11%  "
$ws.Range("E40").Value2 = "  -2.37%  "
$ws.Range("D41").Value2 = "'0.110"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  -1.47%  "
$ws.Range("E42").Value2 = "  -3.03%  "
$ws.Range("D43").Value2 = "'2.000.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value2 = "  -2.48%  "
$ws.Range("D45").Value2 = "'18.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +4.00%  "
$ws.Range("D46").Value2 = "'9.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  -3.53%  "
$ws.Range("D47").Value2 = "'2.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  -7.31%  "
$ws.Range("D48").Value2 = "'2.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -2.35%  "
$ws.Range("D49").Value2 = "'2.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +4.33%  "
$ws.Range("D50").Value2 = "'53.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -1.57%  "
$ws.Range("D51").Value2 = "'2.503.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -1.07%  "
